# Auto-generated Excel COM-interop script to apply updated market/profit data
# to the Seraph_Profits workbook, per the scheduled runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 814.7692
$ws.Range("I28").Value = 849.7
$ws.Range("J28").Value = 698.3333
$ws.Range("K28").Value = 849.7
$ws.Range("L28").Value = 698.3333
$ws.Range("M28").Value = -364.7
$ws.Range("N28").Value = -1668.3333
$ws.Range("H33").Value = 80.583336
$ws.Range("I33").Value = 85.181816
$ws.Range("J33").Value = 30
$ws.Range("K33").Value = 85.181816
$ws.Range("L33").Value = 30
$ws.Range("M33").Value = 143.818184
$ws.Range("N33").Value = -488
$ws.Range("H41").Value = 400
$ws.Range("I41").Value = 400
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 400
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = 40
$ws.Range("N41").ClearContents()
$ws.Range("H53").Value = 227.5
$ws.Range("I53").Value = 132.71428
$ws.Range("J53").Value = 360.2
$ws.Range("K53").Value = 132.71428
$ws.Range("L53").Value = 360.2
$ws.Range("M53").Value = 504.28572
$ws.Range("N53").Value = -1634.2
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()
$ws.Range("H106").Value = 39713.145
$ws.Range("I106").Value = 44665.5
$ws.Range("J106").Value = 9999
$ws.Range("K106").Value = 44665.5
$ws.Range("L106").Value = 9999
$ws.Range("M106").Value = -44034.5
$ws.Range("N106").Value = -11261
$ws.Range("H107").Value = 240.77777
$ws.Range("I107").Value = 224.75
$ws.Range("J107").Value = 369
$ws.Range("K107").Value = 224.75
$ws.Range("L107").Value = 369
$ws.Range("M107").Value = 1695.25
$ws.Range("N107").Value = -4209

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 65000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 65000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 65000
$ws.Range("N64").Value = -65496
$ws.Range("H67").Value = 65000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 65000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 65000
$ws.Range("N67").Value = -66716
$ws.Range("H94").Value = 65000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 65000
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 65000
$ws.Range("N94").Value = -66802
$ws.Range("H110").Value = 7079.9287
$ws.Range("I110").Value = 5318.8184
$ws.Range("J110").Value = 13537.333
$ws.Range("K110").Value = 5318.8184
$ws.Range("L110").Value = 13537.333
$ws.Range("M110").Value = -3273.8184
$ws.Range("N110").Value = -17627.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 50000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 50000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 50000
$ws.Range("N88").Value = -50812
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 50000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 50000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 50000
$ws.Range("N91").Value = -52808
$ws.Range("M91").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H99").Value = 14214.479
$ws.Range("I99").Value = 11823.333
$ws.Range("J99").Value = 15751.643
$ws.Range("K99").Value = 11823.333
$ws.Range("L99").Value = 15751.643
$ws.Range("M99").Value = -10325.333
$ws.Range("N99").Value = -18747.643
$ws.Range("H105").Value = 2077.5
$ws.Range("I105").Value = 1450
$ws.Range("J105").Value = 2705
$ws.Range("K105").Value = 1450
$ws.Range("L105").Value = 2705
$ws.Range("M105").Value = 297
$ws.Range("N105").Value = -6199
$ws.Range("H107").Value = 872.86957
$ws.Range("I107").Value = 490.54544
$ws.Range("J107").Value = 1223.3334
$ws.Range("K107").Value = 490.54544
$ws.Range("L107").Value = 1223.3334
$ws.Range("M107").Value = 1429.45456
$ws.Range("N107").Value = -5063.3334
$ws.Range("H126").Value = 14214.479
$ws.Range("I126").Value = 11823.333
$ws.Range("J126").Value = 15751.643
$ws.Range("K126").Value = 35469.999
$ws.Range("L126").Value = 47254.929
$ws.Range("M126").Value = -32999.999
$ws.Range("N126").Value = -52194.929
$ws.Range("H132").Value = 4286.5
$ws.Range("I132").Value = 3304.75
$ws.Range("J132").Value = 6250
$ws.Range("K132").Value = 9914.25
$ws.Range("L132").Value = 18750
$ws.Range("M132").Value = -7384.25
$ws.Range("N132").Value = -23810
$ws.Range("H134").Value = 2656.4614
$ws.Range("I134").Value = 2059.7778
$ws.Range("J134").Value = 3999
$ws.Range("K134").Value = 6179.3334
$ws.Range("L134").Value = 11997
$ws.Range("M134").Value = -3644.3334
$ws.Range("N134").Value = -17067

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H122").Value = 460.58334
$ws.Range("I122").Value = 339.6
$ws.Range("J122").Value = 547
$ws.Range("K122").Value = 3056.4
$ws.Range("L122").Value = 4923
$ws.Range("M122").Value = -606.4000000000001
$ws.Range("N122").Value = -9823

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 60000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 60000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 60000
$ws.Range("N62").Value = -61372
$ws.Range("H65").Value = 60000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 60000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 180000
$ws.Range("N65").Value = -186864
$ws.Range("H101").Value = 23000
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 23000
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 23000
$ws.Range("N101").Value = -29490
$ws.Range("H113").Value = 3566.2942
$ws.Range("I113").Value = 2236.4
$ws.Range("J113").Value = 4120.4165
$ws.Range("K113").Value = 2236.4
$ws.Range("L113").Value = 4120.4165
$ws.Range("M113").Value = -66.40000000000009
$ws.Range("N113").Value = -8460.416499999999
$ws.Range("H132").Value = 2829.3462
$ws.Range("I132").Value = 1973.2727
$ws.Range("J132").Value = 3457.1333
$ws.Range("K132").Value = 5919.8181
$ws.Range("L132").Value = 10371.3999
$ws.Range("M132").Value = -3389.8181
$ws.Range("N132").Value = -15431.3999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 10499.75
$ws.Range("I16").Value = 13666.333
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 13666.333
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -13496.333
$ws.Range("N16").Value = -1340
$ws.Range("H22").Value = 9647.666999999999
$ws.Range("I22").Value = 949.5
$ws.Range("J22").Value = 13996.75
$ws.Range("K22").Value = 949.5
$ws.Range("L22").Value = 13996.75
$ws.Range("M22").Value = -654.5
$ws.Range("N22").Value = -14586.75
$ws.Range("H27").Value = 9647.666999999999
$ws.Range("I27").Value = 949.5
$ws.Range("J27").Value = 13996.75
$ws.Range("K27").Value = 949.5
$ws.Range("L27").Value = 13996.75
$ws.Range("M27").Value = -842.5
$ws.Range("N27").Value = -14210.75
$ws.Range("H61").Value = 6555.273
$ws.Range("I61").Value = 6815.4443
$ws.Range("J61").Value = 5384.5
$ws.Range("K61").Value = 6815.4443
$ws.Range("L61").Value = 5384.5
$ws.Range("M61").Value = -6613.4443
$ws.Range("N61").Value = -5788.5
$ws.Range("H82").Value = 2311.5
$ws.Range("I82").Value = 2588.7273
$ws.Range("J82").Value = 1549.125
$ws.Range("K82").Value = 2588.7273
$ws.Range("L82").Value = 1549.125
$ws.Range("M82").Value = -2227.7273
$ws.Range("N82").Value = -2271.125
$ws.Range("H85").Value = 2311.5
$ws.Range("I85").Value = 2588.7273
$ws.Range("J85").Value = 1549.125
$ws.Range("K85").Value = 2588.7273
$ws.Range("L85").Value = 1549.125
$ws.Range("M85").Value = -1340.7273
$ws.Range("N85").Value = -4045.125
$ws.Range("H100").Value = 1826.1333
$ws.Range("I100").Value = 1639.2
$ws.Range("J100").Value = 2200
$ws.Range("K100").Value = 1639.2
$ws.Range("L100").Value = 2200
$ws.Range("M100").Value = -1098.2
$ws.Range("N100").Value = -3282
$ws.Range("H113").Value = 6555.273
$ws.Range("I113").Value = 6815.4443
$ws.Range("J113").Value = 5384.5
$ws.Range("K113").Value = 6815.4443
$ws.Range("L113").Value = 5384.5
$ws.Range("M113").Value = -4645.4443
$ws.Range("N113").Value = -9724.5
$ws.Range("H132").Value = 5776.6924
$ws.Range("I132").Value = 5099.4
$ws.Range("J132").Value = 6200
$ws.Range("K132").Value = 15298.2
$ws.Range("L132").Value = 18600
$ws.Range("M132").Value = -12768.2
$ws.Range("N132").Value = -23660

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
